$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.304.70"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.601.23"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.73"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.244"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0607"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.00"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "1.825.34"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "1.597.30"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.73"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "26.295.23"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.47"
$ws.Range("E18").Value = "  +7.73%  "
$ws.Range("D19").Value = "0.0₃0721"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.71"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.95"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "1.449.27"
$ws.Range("E33").Value = "  +6.86%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.569"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.79"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "1.738.47"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.917"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.758"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.76"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.41"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.49"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0951"
$ws.Range("E51").Value = "  -2.79%  "
